# Refresh the cryptocurrency snapshot figures (price + 1h volume change).
# Price cells must stay TEXT (they carry formats like "62.838.42" or
# "0.0000141" that a real number would mangle), so for any new price that
# Excel would otherwise auto-detect as a number we briefly force the cell
# to Text format, assign it, then drop back to the "Normal" style so no
# stray formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.838.42"
$ws.Range("E2").Value = "  -1.59%  "

# Row 3
$ws.Range("D3").Value = "2.533.45"
$ws.Range("E3").Value = "  -0.21%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.29%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("E8").Value = "  -1.34%  "

# Row 9
$ws.Range("D9").Value = "2.530.71"
$ws.Range("E9").Value = "  -0.42%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.69%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.18%  "

# Row 12
$ws.Range("E12").Value = "  -0.54%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.97%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.44%  "

# Row 15
$ws.Range("D15").Value = "2.985.34"
$ws.Range("E15").Value = "  -0.22%  "

# Row 16
$ws.Range("D16").Value = "62.760.23"
$ws.Range("E16").Value = "  -1.34%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.99%  "

# Row 18
$ws.Range("D18").Value = "2.532.40"
$ws.Range("E18").Value = "  -0.13%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.97%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "332.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.17%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.62%  "

# Row 23
$ws.Range("E23").Value = "  -0.61%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.45%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.22%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.24%  "

# Row 27
$ws.Range("E27").Value = "  -0.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.52%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.70%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0802"
$ws.Range("E31").Value = "  -2.80%  "

# Row 32
$ws.Range("E32").Value = "  -2.07%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.06%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "395.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.27%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.38%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.396"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.29%  "

# Row 38
$ws.Range("E38").Value = "  +0.01%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.98%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.07%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.14%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.28%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "150.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.56%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.94%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.79%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0527"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.68%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.596"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.40%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0960"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.64%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0236"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.04%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.45%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
